$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume(1h) (E) columns are treated as text,
# matching the original inline-string cell contents (e.g. "216.35", "1.800"),
# so Excel does not auto-convert numeric-looking values into real numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.213.65"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -4.11%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.655.02"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -3.53%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.35"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -3.66%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5138"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.00%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.007"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.88%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06456"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -3.46%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.94"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -4.56%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07848"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.11%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.305"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -4.18%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.654.91"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.53%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.883.56"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5528"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -4.61%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8040"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.88%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.19"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -5.25%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.217.18"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -4.19%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.00%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "210.74"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -4.93%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.419"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -4.94%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.10"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -3.29%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.024"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.01%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.02%  "

$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.84"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.44%  "

$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.800"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +5.54%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1176"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.52%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.28%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.86"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.27%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05105"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -5.09%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -4.15%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.367"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.26%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.234"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -5.28%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.562"

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.735"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -4.19%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9257"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.68%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.352"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.84%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5734"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.80%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.165.86"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.65%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01591"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -3.71%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.88%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.713"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.13%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8242"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.86%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "100.37"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.797.76"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.27%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.99%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4548"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.68%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "55.45"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -4.08%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.006"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.893"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.76%  "
